# Add "Program" column to the Init sheet picks/spawning parser template,
# matching the new Program column already present on the Picking sheet,
# and make the Init sheet the active tab instead of Picking.

$wb = $excel.ActiveWorkbook

$wsInit = $wb.Worksheets.Item("Init")

# Grab the text of the comments that currently live on the "Tray" (G3) and
# "Fecundity" (H3) header cells, before anything shifts, so they can be
# recreated at their new locations after the column insert below.
$mandatoryCommentText = $wsInit.Range("G3").Comment.Text()
$fecundityCommentText = $wsInit.Range("H3").Comment.Text()

# --- Init sheet: insert a new "Program" column before the old "Cross" column (E) ---
$wsInit.Columns("E:E").Insert()
$wsInit.Range("E3").Value = "Program"

# Remove the old sample data rows (4:6) that are no longer needed.
$wsInit.Rows("4:6").Delete()

# Comments are anchored to cells and are not shifted by the column insert, so
# remove the now-stale comments left behind at G3/H3 and recreate them at
# their new homes (H3/I3), and add a new comment describing the Program column.
$wsInit.Range("G3").Comment.Delete()
$wsInit.Range("H3").Comment.Delete()

$wsInit.Range("H3").AddComment($mandatoryCommentText)
$wsInit.Range("I3").AddComment($fecundityCommentText)
$wsInit.Range("E3").AddComment("Program of pairing set in spawning. Optional, must match code in db.")

# Make "Init" the active/selected sheet and cell, instead of "Picking".
$wsInit.Range("E3").Select()
$wsInit.Activate()
